$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Actual Consumption values (column A, rows 2-53)
$aVals = @(6690,6607,6577,6544,6559,6500,6453,6470,6389,6405,6407,6368,6439,6382,6387,6437,6415,6441,6508,6505,6668,6737,6758,6903,7086,7258,7369,7663,7974,8225,8384,8480,8709,8816,8841,8846,8762,8758,8685,8616,8533,8358,8257,8144,8012,7945,7824,7727,7776,7686,7704,7659)

# New Timestamp serial values (column B, rows 2-53)
$bVals = @(46055.95833333334,46055.96875,46055.97916666666,46055.98958333334,46056.0,46056.01041666666,46056.02083333334,46056.03125,46056.04166666666,46056.05208333334,46056.0625,46056.07291666666,46056.08333333334,46056.09375,46056.10416666666,46056.11458333334,46056.125,46056.13541666666,46056.14583333334,46056.15625,46056.16666666666,46056.17708333334,46056.1875,46056.19791666666,46056.20833333334,46056.21875,46056.22916666666,46056.23958333334,46056.25,46056.26041666666,46056.27083333334,46056.28125,46056.29166666666,46056.30208333334,46056.3125,46056.32291666666,46056.33333333334,46056.34375,46056.35416666666,46056.36458333334,46056.375,46056.38541666666,46056.39583333334,46056.40625,46056.41666666666,46056.42708333334,46056.4375,46056.44791666666,46056.45833333334,46056.46875,46056.47916666666,46056.48958333334)

for ($i = 0; $i -lt $aVals.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}

# Ensure the new rows (50-53) inherit the Timestamp column's date/time number format
$ws.Range("B50:B53").NumberFormat = "YYYY-MM-DD HH:MM:SS"
